$wb = $excel.ActiveWorkbook

# Update 2025 (and, where the daily refresh nudged a late-reported
# incident, 2024) violent-crime counts for the 2025-10-05 data pull.

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5164  # was 5134
$ws.Range("L3").Value = 5559  # was 5532
$ws.Range("L4").Value = 1354  # was 1350
$ws.Range("L5").Value = 332  # was 330
$ws.Range("L6").Value = 4628  # was 4608
$ws.Range("L7").Value = 17037  # was 16954

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 61  # was 60
$ws.Range("L6").Value = 129  # was 128
$ws.Range("L8").Value = 1123  # was 1117
$ws.Range("L10").Value = 109  # was 108
$ws.Range("L12").Value = 39  # was 38
$ws.Range("L15").Value = 129  # was 128
$ws.Range("L19").Value = 461  # was 460
$ws.Range("L20").Value = 420  # was 419
$ws.Range("L23").Value = 187  # was 184
$ws.Range("L25").Value = 102  # was 101
$ws.Range("L29").Value = 940  # was 934
$ws.Range("L33").Value = 789  # was 784
$ws.Range("L34").Value = 101  # was 99
$ws.Range("L36").Value = 217  # was 216
$ws.Range("L37").Value = 649  # was 643
$ws.Range("L40").Value = 49  # was 48
$ws.Range("K42").Value = 1030  # was 1029
$ws.Range("L42").Value = 557  # was 555
$ws.Range("L43").Value = 124  # was 123
$ws.Range("L45").Value = 31  # was 30
$ws.Range("L47").Value = 117  # was 116
$ws.Range("L51").Value = 212  # was 210
$ws.Range("L52").Value = 344  # was 341
$ws.Range("L54").Value = 365  # was 363
$ws.Range("L55").Value = 174  # was 173
$ws.Range("L59").Value = 30  # was 31
$ws.Range("K63").Value = 173  # was 174
$ws.Range("L64").Value = 114  # was 113
$ws.Range("L65").Value = 330  # was 329
$ws.Range("L66").Value = 45  # was 44
$ws.Range("L67").Value = 587  # was 583
$ws.Range("L71").Value = 47  # was 46
$ws.Range("L73").Value = 136  # was 134
$ws.Range("L76").Value = 262  # was 260
$ws.Range("L77").Value = 112  # was 110
$ws.Range("L78").Value = 217  # was 215
$ws.Range("L79").Value = 454  # was 448
$ws.Range("L84").Value = 167  # was 164
$ws.Range("L85").Value = 861  # was 858
$ws.Range("L86").Value = 120  # was 119
$ws.Range("L94").Value = 213  # was 211
$ws.Range("L99").Value = 302  # was 300
$ws.Range("L101").Value = 17037  # was 16954

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 260  # was 259
$ws.Range("L3").Value = 352  # was 350
$ws.Range("L7").Value = 861  # was 858

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 111  # was 109
$ws.Range("L3").Value = 111  # was 110
$ws.Range("L7").Value = 344  # was 341

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 328  # was 325
$ws.Range("L3").Value = 385  # was 383
$ws.Range("L6").Value = 289  # was 288
$ws.Range("L7").Value = 1123  # was 1117

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 218  # was 216
$ws.Range("L3").Value = 275  # was 273
$ws.Range("L5").Value = 19  # was 18
$ws.Range("L7").Value = 789  # was 784

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 196  # was 195
$ws.Range("L3").Value = 221  # was 217
$ws.Range("L6").Value = 178  # was 177
$ws.Range("L7").Value = 649  # was 643

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 121  # was 120
$ws.Range("L7").Value = 330  # was 329

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 125  # was 124
$ws.Range("L6").Value = 64  # was 63
$ws.Range("L7").Value = 302  # was 300

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 227  # was 225
$ws.Range("L5").Value = 16  # was 15
$ws.Range("L6").Value = 134  # was 133
$ws.Range("L7").Value = 587  # was 583

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 54  # was 53
$ws.Range("L3").Value = 58  # was 57
$ws.Range("L4").Value = 6  # was 5
$ws.Range("L7").Value = 167  # was 164

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 67  # was 66
$ws.Range("L6").Value = 178  # was 177
$ws.Range("L7").Value = 365  # was 363

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 290  # was 286
$ws.Range("L3").Value = 352  # was 351
$ws.Range("L6").Value = 234  # was 233
$ws.Range("L7").Value = 940  # was 934

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 143  # was 142
$ws.Range("L7").Value = 461  # was 460

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 48  # was 47
$ws.Range("L6").Value = 122  # was 121
$ws.Range("L7").Value = 262  # was 260

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L6").Value = 25  # was 24
$ws.Range("L7").Value = 129  # was 128

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 155  # was 154
$ws.Range("K4").Value = 48  # was 47
$ws.Range("L6").Value = 153  # was 152
$ws.Range("K7").Value = 1030  # was 1029
$ws.Range("L7").Value = 557  # was 555

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L6").Value = 31  # was 30
$ws.Range("L7").Value = 109  # was 108

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 60  # was 58
$ws.Range("L7").Value = 217  # was 215

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 53  # was 52
$ws.Range("L7").Value = 174  # was 173

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 48  # was 46
$ws.Range("L6").Value = 49  # was 48
$ws.Range("L7").Value = 187  # was 184

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L6").Value = 105  # was 99
$ws.Range("L7").Value = 454  # was 448

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 32  # was 31
$ws.Range("L7").Value = 114  # was 113

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 130  # was 129
$ws.Range("L7").Value = 420  # was 419

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L4").Value = 18  # was 17
$ws.Range("L7").Value = 217  # was 216

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L2").Value = 33  # was 31
$ws.Range("L7").Value = 101  # was 99

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L3").Value = 50  # was 49
$ws.Range("L6").Value = 81  # was 80
$ws.Range("L7").Value = 213  # was 211

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L2").Value = 36  # was 35
$ws.Range("L7").Value = 102  # was 101

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 38  # was 37
$ws.Range("L7").Value = 117  # was 116

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 46  # was 45
$ws.Range("L7").Value = 129  # was 128

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L3").Value = 13  # was 12
$ws.Range("L7").Value = 45  # was 44

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L3").Value = 40  # was 38
$ws.Range("L7").Value = 136  # was 134

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L6").Value = 7  # was 8
$ws.Range("L7").Value = 30  # was 31

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L2").Value = 14  # was 13
$ws.Range("L7").Value = 61  # was 60

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 65  # was 64
$ws.Range("L7").Value = 120  # was 119

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 67  # was 66
$ws.Range("L4").Value = 31  # was 30
$ws.Range("L7").Value = 212  # was 210

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L3").Value = 39  # was 38
$ws.Range("L7").Value = 124  # was 123

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L2").Value = 21  # was 20
$ws.Range("L7").Value = 47  # was 46

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L2").Value = 38  # was 37
$ws.Range("L6").Value = 25  # was 24
$ws.Range("L7").Value = 112  # was 110

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("L2").Value = 6  # was 5
$ws.Range("L7").Value = 31  # was 30

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L3").Value = 23  # was 22
$ws.Range("L7").Value = 49  # was 48

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("L6").Value = 10  # was 9
$ws.Range("L7").Value = 39  # was 38
